$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Refreshed auth tokens / ids produced by re-running the user-data test
# fixtures (setUp) for daniel5f, Jorge2525 and mario35.
$ws.Range("C2").Value = "eyJhbGciOiJIUzI1NiIsInR5cCI6IkpXVCJ9.eyJ1c2VyTmFtZSI6ImRhbmllbDVmIiwicGFzc3dvcmQiOiJBejI1Mjg4QCIsImlhdCI6MTcwMjc3Njk0NH0.8byoKoaKXjfBNXXvEsH5OkGmF8PKd5ag66N5BdeWKCE"
$ws.Range("D2").Value = "d9389dba-2a4a-4b86-8644-99c86ea65fb8"

$ws.Range("C3").Value = "eyJhbGciOiJIUzI1NiIsInR5cCI6IkpXVCJ9.eyJ1c2VyTmFtZSI6IkpvcmdlMjUyNSIsInBhc3N3b3JkIjoiYXNUMzU2NDQ0QCIsImlhdCI6MTcwMjc3Njk0Nn0.whkZ-sKpdfPyNB-g-ZadoKc_au-FYmC4vgDws_Oajdk"
$ws.Range("D3").Value = "bf2d2bf4-e6ca-4e2b-a5fd-276c072ed84d"

$ws.Range("C4").Value = "eyJhbGciOiJIUzI1NiIsInR5cCI6IkpXVCJ9.eyJ1c2VyTmFtZSI6Im1hcmlvMzUiLCJwYXNzd29yZCI6Im1BcmlvdXVnQDMiLCJpYXQiOjE3MDI3NzY5NDd9.9Hr6MClQTOSdaKWpR6oyXu0LmB_crjc-cVTACNnydo8"
$ws.Range("D4").Value = "3632eb03-7ec4-4ade-8cf8-6995e8131469"
